$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 127; this shifts existing rows 127-192 down to 128-193
# and Excel extends the used range / dimension automatically.
$ws.Rows.Item(127).Insert()

# Populate the newly inserted row 127 with the new weekly record.
$ws.Cells.Item(127, 1).Value = 3
$ws.Cells.Item(127, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(127, 3).Value = "Coquimbo"
$ws.Cells.Item(127, 4).Value = 44455
$ws.Cells.Item(127, 5).Value = 5
$ws.Cells.Item(127, 6).Value = 100112009
$ws.Cells.Item(127, 7).Value = "Acelga"
$ws.Cells.Item(127, 8).Value = "Sin especificar"
$ws.Cells.Item(127, 9).Value = "Primera"
$ws.Cells.Item(127, 10).Value = 250
$ws.Cells.Item(127, 11).Value = 2000
$ws.Cells.Item(127, 12).Value = 2200
$ws.Cells.Item(127, 13).Value = 2104
$ws.Cells.Item(127, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(127, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(127, 16).Value = 351
$ws.Cells.Item(127, 17).Value = 6
$ws.Cells.Item(127, 18).Value = "Hortaliza"
